# "Added charge per appt" - add a new "ChargePerSquareFootage" column (E) to
# the apartment details sheet, and rename the 4th row's "Owner" role to
# "Tenant" (C4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header in E1
$ws.Range("E1").Value = "ChargePerSquareFootage"

# Charge-per-square-footage values for each apartment row
$ws.Range("E2").Value = 3.25
$ws.Range("E3").Value = 3.25
$ws.Range("E4").Value = 3.58

# Row 4 (SRK) is a Tenant rather than an Owner
$ws.Range("C4").Value = "Tenant"

# Size the new column to fit its (longer) header text, same as the other
# bestFit columns in the sheet
$ws.Columns.Item(5).ColumnWidth = 20

# Match the selection left behind after entering the last new value
$ws.Range("E4").Select()
